$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1795.4762
$ws.Range("I98").Value = 1174.9375
$ws.Range("J98").Value = 3781.2
$ws.Range("K98").Value = 1174.9375
$ws.Range("L98").Value = 3781.2
$ws.Range("M98").Value = 323.0625
$ws.Range("N98").Value = -6777.2
# Row 111
$ws.Range("H111").Value = 2129.0435
$ws.Range("I111").Value = 1544.0625
$ws.Range("J111").Value = 3466.1428
$ws.Range("K111").Value = 4632.1875
$ws.Range("L111").Value = 10398.4284
$ws.Range("M111").Value = -1565.1875
$ws.Range("N111").Value = -16532.4284
# Row 113
$ws.Range("H113").Value = 4425.727
$ws.Range("I113").Value = 3478
$ws.Range("J113").Value = 5215.5
$ws.Range("K113").Value = 3478
$ws.Range("L113").Value = 5215.5
$ws.Range("M113").Value = -224
$ws.Range("N113").Value = -11723.5
# Row 115
$ws.Range("H115").Value = 1761.6666
$ws.Range("I115").Value = 1193.3334
$ws.Range("J115").Value = 3466.6667
$ws.Range("K115").Value = 3580.0002
$ws.Range("L115").Value = 10400.0001
$ws.Range("M115").Value = -2013.0002
$ws.Range("N115").Value = -13534.0001
# Row 122
$ws.Range("H122").Value = 1795.4762
$ws.Range("I122").Value = 1174.9375
$ws.Range("J122").Value = 3781.2
$ws.Range("K122").Value = 3524.8125
$ws.Range("L122").Value = 11343.6
$ws.Range("M122").Value = -1074.8125
$ws.Range("N122").Value = -16243.6

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 1845.9667
$ws.Range("I122").Value = 1584.8182
$ws.Range("K122").Value = 4754.4546
$ws.Range("M122").Value = -2304.4546

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3958.9268
$ws.Range("I134").Value = 3560.1875
$ws.Range("K134").Value = 10680.5625
$ws.Range("M134").Value = -8145.5625

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 627.63635
$ws.Range("I10").Value = 600.5
$ws.Range("J10").Value = 700
$ws.Range("K10").Value = 600.5
$ws.Range("L10").Value = 700
$ws.Range("M10").Value = -461.5
$ws.Range("N10").Value = -978
# Row 16
$ws.Range("H16").Value = 1028.7142
$ws.Range("I16").Value = 1000.2
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 1000.2
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -713.2
$ws.Range("N16").Value = -1674
# Row 74
$ws.Range("H74").Value = 25018.75
$ws.Range("J74").Value = 25018.75
$ws.Range("L74").Value = 25018.75
$ws.Range("N74").Value = -26766.75
# Row 77
$ws.Range("H77").Value = 25018.75
$ws.Range("J77").Value = 25018.75
$ws.Range("L77").Value = 75056.25
$ws.Range("N77").Value = -83792.25
# Row 113
$ws.Range("H113").Value = 1028.7142
$ws.Range("I113").Value = 1000.2
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1000.2
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1169.8
$ws.Range("N113").Value = -5440
# Row 122
$ws.Range("H122").Value = 945.75
$ws.Range("I122").Value = 879.73914
$ws.Range("J122").Value = 1114.4445
$ws.Range("K122").Value = 2639.21742
$ws.Range("L122").Value = 3343.3335
$ws.Range("M122").Value = -189.2174199999999
$ws.Range("N122").Value = -8243.333500000001
# Row 132
$ws.Range("H132").Value = 2356.7334
$ws.Range("I132").Value = 2055.3333
$ws.Range("J132").Value = 4315.8335
$ws.Range("K132").Value = 6165.999899999999
$ws.Range("L132").Value = 12947.5005
$ws.Range("M132").Value = -3635.999899999999
$ws.Range("N132").Value = -18007.5005

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 12542.625
$ws.Range("I10").Value = 12542.625
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 37627.875
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -37488.875
$ws.Range("N10").ClearContents()
# Row 68
$ws.Range("H68").Value = 1083.1666
$ws.Range("I68").Value = 1083.1666
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3249.4998
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2438.4998
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 1083.1666
$ws.Range("I71").Value = 1083.1666
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9748.499400000001
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5692.499400000001
$ws.Range("N71").ClearContents()
# Row 113
$ws.Range("H113").Value = 565.3158
$ws.Range("I113").Value = 567.5714
$ws.Range("J113").Value = 564
$ws.Range("K113").Value = 1702.7142
$ws.Range("L113").Value = 1692
$ws.Range("M113").Value = 467.2857999999999
$ws.Range("N113").Value = -6032
# Row 131
$ws.Range("H131").Value = 896.91
$ws.Range("J131").Value = 897.0707
$ws.Range("L131").Value = 2691.2121
$ws.Range("N131").Value = -12771.2121

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 6416.25
$ws.Range("I12").Value = 5000
$ws.Range("J12").Value = 6888.3335
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6888.3335
$ws.Range("M12").Value = -4860
$ws.Range("N12").Value = -7168.3335
# Row 98
$ws.Range("H98").Value = 25933.334
$ws.Range("J98").Value = 25933.334
$ws.Range("L98").Value = 25933.334
$ws.Range("N98").Value = -31923.334
# Row 102
$ws.Range("H102").Value = 1131.6923
$ws.Range("I102").Value = 1104.174
$ws.Range("J102").Value = 1342.6666
$ws.Range("K102").Value = 1104.174
$ws.Range("L102").Value = 1342.6666
$ws.Range("M102").Value = 517.826
$ws.Range("N102").Value = -4586.6666
# Row 113
$ws.Range("H113").Value = 1840.4445
$ws.Range("I113").Value = 1274.5714
$ws.Range("J113").Value = 2200.5454
$ws.Range("K113").Value = 1274.5714
$ws.Range("L113").Value = 2200.5454
$ws.Range("M113").Value = 895.4286
$ws.Range("N113").Value = -6540.5454
# Row 122
$ws.Range("H122").Value = 53861.74
$ws.Range("I122").Value = 60127.94
$ws.Range("J122").Value = 599
$ws.Range("K122").Value = 180383.82
$ws.Range("L122").Value = 1797
$ws.Range("M122").Value = -177933.82
$ws.Range("N122").Value = -6697

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2600.077
$ws.Range("I7").Value = 1367.75
$ws.Range("J7").Value = 3147.7778
$ws.Range("K7").Value = 1367.75
$ws.Range("L7").Value = 3147.7778
$ws.Range("M7").Value = -1255.75
$ws.Range("N7").Value = -3371.7778
# Row 40
$ws.Range("H40").Value = 4468.615
$ws.Range("I40").Value = 4657.6665
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 4657.6665
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -4521.6665
$ws.Range("N40").Value = -2472
# Row 46
$ws.Range("H46").Value = 896.2917
$ws.Range("I46").Value = 1200.2
$ws.Range("J46").Value = 816.3158
$ws.Range("K46").Value = 1200.2
$ws.Range("L46").Value = 816.3158
$ws.Range("M46").Value = -1012.2
$ws.Range("N46").Value = -1192.3158
# Row 122
$ws.Range("H122").Value = 2769.1714
$ws.Range("I122").Value = 2247.125
$ws.Range("J122").Value = 3208.7896
$ws.Range("K122").Value = 6741.375
$ws.Range("L122").Value = 9626.3688
$ws.Range("M122").Value = -4291.375
$ws.Range("N122").Value = -14526.3688
# Row 126
$ws.Range("H126").Value = 2600.077
$ws.Range("I126").Value = 1367.75
$ws.Range("J126").Value = 3147.7778
$ws.Range("K126").Value = 4103.25
$ws.Range("L126").Value = 9443.3334
$ws.Range("M126").Value = -1633.25
$ws.Range("N126").Value = -14383.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 112
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
# Row 126
$ws.Range("H126").Value = 1376.2609
$ws.Range("I126").Value = 829
$ws.Range("J126").Value = 2402.375
$ws.Range("K126").Value = 2487
$ws.Range("L126").Value = 7207.125
$ws.Range("M126").Value = -17
$ws.Range("N126").Value = -12147.125
